# Update countries & provincias Spain
# Applies the "data refresh" edit described in the commit:
#  - Updates the "Datos actualizados..." timestamp string
#  - Updates case numbers for several countries (simple value refresh)
#  - Uruguay overtakes Afganistan in the ranking (rows 88/89 swap identity)
#  - San Marino overtakes Niger in the ranking (rows 101/102 swap identity)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the timestamp row (A1) ---
$ws.Range("A1").Value = "Datos actualizados a 8 de Abril de 2020 a las 01:22"

# --- Simple numeric refreshes (country keeps its row/rank) ---

# Row 4: Estados Unidos
$ws.Range("B4").Value = 395277
$ws.Range("C4").Value = 28273
$ws.Range("E4").Value = 360819
$ws.Range("G4").Value = 1913
$ws.Range("H4").Value = 12784

# Row 16: Canada
$ws.Range("B16").Value = 17897
$ws.Range("C16").Value = 1230
$ws.Range("E16").Value = 13488
$ws.Range("G16").Value = 58
$ws.Range("H16").Value = 381

# Row 17: Brasil
$ws.Range("B17").Value = 14034
$ws.Range("C17").Value = 1851
$ws.Range("E17").Value = 13221

# Row 24: Noruega
$ws.Range("B24").Value = 6086
$ws.Range("C24").Value = 221
$ws.Range("E24").Value = 5965

# Row 85: Uzbekistan
$ws.Range("B85").Value = 520
$ws.Range("C85").Value = 63
$ws.Range("E85").Value = 488

# Row 117: Venezuela
$ws.Range("B117").Value = 166
$ws.Range("C117").Value = 1
$ws.Range("E117").Value = 94

# Row 139: Jamaica
$ws.Range("D139").Value = 9
$ws.Range("E139").Value = 51

# --- Ranking swaps ---

# Uruguay's updated totals now exceed Afganistan's, so Uruguay moves up to
# row 88 (taking the old Afganistan row's place) and Afganistan drops to
# row 89 with the data Uruguay/Afganistan previously held at row 88.
$ws.Range("A88").Value = "Uruguay"
$ws.Range("B88").Value = 424
$ws.Range("C88").Value = 9
$ws.Range("D88").Value = 150
$ws.Range("E88").Value = 267
$ws.Range("F88").Value = 14
$ws.Range("G88").Value = 1
$ws.Range("H88").Value = 7

$ws.Range("A89").Value = "Afganistan"
$ws.Range("B89").Value = 423
$ws.Range("C89").Value = 56
$ws.Range("D89").Value = 18
$ws.Range("E89").Value = 391
$ws.Range("F89").Value = 0
$ws.Range("G89").Value = 3
$ws.Range("H89").Value = 14

# San Marino's updated totals now exceed Niger's, so San Marino moves up to
# row 101 (taking the old Niger row's place) and Niger drops to row 102.
$ws.Range("A101").Value = "San Marino"
$ws.Range("B101").Value = 279
$ws.Range("C101").Value = 2
$ws.Range("D101").Value = 40
$ws.Range("E101").Value = 205
$ws.Range("F101").Value = 14
$ws.Range("G101").Value = 2
$ws.Range("H101").Value = 34

$ws.Range("A102").Value = "Niger"
$ws.Range("B102").Value = 278
$ws.Range("C102").Value = 25
$ws.Range("D102").Value = 26
$ws.Range("E102").Value = 241
$ws.Range("F102").Value = 0
$ws.Range("G102").Value = 1
$ws.Range("H102").Value = 11
